# Builds a string consisting of $n copies of $ch (PowerShell's "*" operator
# for string repetition is not supported by this interpreter).
function Repeat($ch, $n) {
    $s = ""
    for ($i = 0; $i -lt $n; $i++) {
        $s = $s + $ch
    }
    return $s
}

$d = $word.ActiveDocument
$tab = [string][char]9

# ---------------------------------------------------------------------------
# Paragraph 1 (signature block just above "CSC President"):
#   BEFORE: ${cscAdviser}<5 tabs><45 spaces>${oicOsa}
#   AFTER:  ${cscPresident}<5 tabs><13 spaces>${cscAdviser}
# ---------------------------------------------------------------------------

# Step 1: rename the first placeholder token cscAdviser -> cscPresident.
# At this point in the document "cscAdviser" is unique, so this is safe.
$find1 = $d.Content.Find
$r1 = $find1.Execute("cscAdviser", $true, $false, $false, $false, $false, $true, 1, $false, "cscPresident", 2)

# Step 2: rename the second placeholder token oicOsa -> cscAdviser.
$find2 = $d.Content.Find
$r2 = $find2.Execute("oicOsa", $true, $false, $false, $false, $false, $true, 1, $false, "cscAdviser", 2)

# Step 3: shrink the run of spaces between the tabs and "${" from 45 to 13
# (search text has no tabs in it, so the match stays inside the existing
# text run and the surrounding <w:tab/> runs are left untouched).
$old3 = Repeat " " 45
$new3 = Repeat " " 13
$find3 = $d.Content.Find
$r3 = $find3.Execute($old3, $true, $false, $false, $false, $false, $true, 1, $false, $new3, 2)

# ---------------------------------------------------------------------------
# Paragraph 2 ("CSC President ... Adviser" signature line):
#   BEFORE: CSC President<7 tabs><20 spaces>Adviser<5 tabs>
#   AFTER:  CSC President<7 tabs><1 space>Adviser<11 tabs>
# ---------------------------------------------------------------------------

# Step 4: shrink the spaces immediately before "Adviser" from 20 down to 1
# (again, no tabs in the search text, so this is a plain text substitution
# inside the existing run).
$old4 = (Repeat " " 20) + "Adviser"
$new4 = " Adviser"
$find4 = $d.Content.Find
$r4 = $find4.Execute($old4, $true, $false, $false, $false, $false, $true, 1, $false, $new4, 2)

# Step 5: grow the trailing run of tabs from 5 to 11 (6 extra tabs), by
# inserting right at the boundary just after the last existing tab
# character so the new text picks up that run's formatting.
$p2 = $d.Paragraphs.Item(66)
$rng2 = $p2.Range
$lastTabPos = $rng2.End - 2
$insertPoint = $d.Range($lastTabPos + 1, $lastTabPos + 1)
$insertPoint.InsertBefore((Repeat $tab 6))

Write-Output "r1=$r1 r2=$r2 r3=$r3 r4=$r4"
